# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> used by the slide master (was the "Integral" theme)
#   ppt/theme/theme2.xml -> used by the notes master (was the "Office Theme")
# The author's edit swaps which theme each master points at: the slide
# master now carries the stock "Office Theme" color palette (previously
# "Integral"), while the notes master keeps the "Integral" colors it
# inherited from the slide master before. The font scheme and format
# scheme (fills/lines/effects) are identical between the two themes
# already, so the only real content change is the 12-color scheme
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# Apply it by pushing the standard Office theme palette onto the
# presentation's live theme color scheme via the DrawingML
# ThemeColorScheme object (2007+ object model, 12 slots, in the same
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order as OOXML's <a:clrScheme>).

function ConvertTo-BgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # OLE/VBA RGB colors are packed as 0x00BBGGRR.
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target "Office Theme" color scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i + 1).RGB = ConvertTo-BgrInt $officeThemeColors[$i]
}
